# Updates cryptos list prices (column D) and volume-change percentages
# (column E) for rows 2-51, mirroring a scheduled GitHub Actions data refresh.
# Values are plain text in the source workbook (e.g. "332.49", "  -1.08%  "),
# so each cell is forced to Text format before the write (and its original
# style/format is restored immediately after) to stop Excel's COM layer from
# auto-coercing number-looking strings ("332.49", "1.002", ...) into numerics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.565.92"
Set-TextValue $ws.Range("E2") "  -1.43%  "
Set-TextValue $ws.Range("D3") "1.846.22"
Set-TextValue $ws.Range("E3") "  -2.31%  "
Set-TextValue $ws.Range("E4") "  -1.18%  "
Set-TextValue $ws.Range("D5") "332.49"
Set-TextValue $ws.Range("E5") "  -1.08%  "
Set-TextValue $ws.Range("E6") "  -1.24%  "
Set-TextValue $ws.Range("D7") "0.4633"
Set-TextValue $ws.Range("E7") "  -1.70%  "
Set-TextValue $ws.Range("D8") "0.3851"
Set-TextValue $ws.Range("E8") "  -2.36%  "
Set-TextValue $ws.Range("D9") "45.93"
Set-TextValue $ws.Range("E9") "  -1.90%  "
Set-TextValue $ws.Range("D10") "0.07895"
Set-TextValue $ws.Range("E10") "  -1.36%  "
Set-TextValue $ws.Range("D11") "0.9923"
Set-TextValue $ws.Range("E11") "  -2.57%  "
Set-TextValue $ws.Range("D12") "21.45"
Set-TextValue $ws.Range("E12") "  -1.48%  "
Set-TextValue $ws.Range("D13") "1.856.69"
Set-TextValue $ws.Range("E13") "  -2.02%  "
Set-TextValue $ws.Range("D14") "5.908"
Set-TextValue $ws.Range("E14") "  -1.60%  "
Set-TextValue $ws.Range("D15") "7.090"
Set-TextValue $ws.Range("E15") "  -1.30%  "
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  -1.53%  "
Set-TextValue $ws.Range("D17") "88.72"
Set-TextValue $ws.Range("E17") "  +0.64%  "
Set-TextValue $ws.Range("D18") "0.06652"
Set-TextValue $ws.Range("E18") "  -1.86%  "
Set-TextValue $ws.Range("D19") "0.00001034"
Set-TextValue $ws.Range("E19") "  -1.70%  "
Set-TextValue $ws.Range("D20") "17.03"
Set-TextValue $ws.Range("E20") "  -0.71%  "
Set-TextValue $ws.Range("E21") "  -1.16%  "
Set-TextValue $ws.Range("D22") "27.595.22"
Set-TextValue $ws.Range("E22") "  -1.37%  "
Set-TextValue $ws.Range("D23") "5.374"
Set-TextValue $ws.Range("E23") "  -2.34%  "
Set-TextValue $ws.Range("E24") "  -0.70%  "
Set-TextValue $ws.Range("D25") "2.304"
Set-TextValue $ws.Range("E25") "  -2.49%  "
Set-TextValue $ws.Range("D26") "158.11"
Set-TextValue $ws.Range("E26") "  -0.90%  "
Set-TextValue $ws.Range("D27") "19.52"
Set-TextValue $ws.Range("E27") "  -2.57%  "
Set-TextValue $ws.Range("D28") "2.095"
Set-TextValue $ws.Range("E28") "  -0.57%  "
Set-TextValue $ws.Range("D29") "5.392"
Set-TextValue $ws.Range("E29") "  -2.00%  "
Set-TextValue $ws.Range("D30") "119.64"
Set-TextValue $ws.Range("E30") "  -1.71%  "
Set-TextValue $ws.Range("D31") "0.9726"
Set-TextValue $ws.Range("E31") "  +0.79%  "
Set-TextValue $ws.Range("D32") "0.09387"
Set-TextValue $ws.Range("E32") "  -1.92%  "
Set-TextValue $ws.Range("D33") "3.576"
Set-TextValue $ws.Range("E33") "  -1.98%  "
Set-TextValue $ws.Range("D34") "5.277"
Set-TextValue $ws.Range("E34") "  -1.58%  "
Set-TextValue $ws.Range("D35") "1.342"
Set-TextValue $ws.Range("E35") "  -1.51%  "
Set-TextValue $ws.Range("D36") "0.06000"
Set-TextValue $ws.Range("E36") "  -2.05%  "
Set-TextValue $ws.Range("D37") "0.02220"
Set-TextValue $ws.Range("E37") "  -1.30%  "
Set-TextValue $ws.Range("D38") "8.284"
Set-TextValue $ws.Range("E38") "  +0.67%  "
Set-TextValue $ws.Range("D39") "1.178"
Set-TextValue $ws.Range("E39") "  -3.01%  "
Set-TextValue $ws.Range("D40") "0.5893"
Set-TextValue $ws.Range("E40") "  -1.42%  "
Set-TextValue $ws.Range("D41") "0.1859"
Set-TextValue $ws.Range("E41") "  -2.56%  "
Set-TextValue $ws.Range("D42") "10.27"
Set-TextValue $ws.Range("E42") "  -0.90%  "
Set-TextValue $ws.Range("D43") "1.243"
Set-TextValue $ws.Range("E43") "  -2.12%  "
Set-TextValue $ws.Range("D44") "0.5576"
Set-TextValue $ws.Range("E44") "  -2.18%  "
Set-TextValue $ws.Range("D45") "12.15"
Set-TextValue $ws.Range("E45") "  -0.75%  "
Set-TextValue $ws.Range("D46") "1.891"
Set-TextValue $ws.Range("E46") "  -2.85%  "
Set-TextValue $ws.Range("D47") "0.06685"
Set-TextValue $ws.Range("E47") "  -2.73%  "
Set-TextValue $ws.Range("D48") "110.50"
Set-TextValue $ws.Range("E48") "  -2.60%  "
Set-TextValue $ws.Range("D49") "1.053"
Set-TextValue $ws.Range("E49") "  -1.69%  "
Set-TextValue $ws.Range("D50") "1.002"
Set-TextValue $ws.Range("E50") "  -1.43%  "
Set-TextValue $ws.Range("D51") "70.05"
Set-TextValue $ws.Range("E51") "  -1.75%  "
